# Insert two new "Buy" trade rows above the existing Feb-10 trade row on the
# "Trading History" sheet, pushing the existing row (date 46063 / 2026-02-10)
# down from row 5 to row 7, and populating the two newly inserted rows
# (5 and 6) with Feb-09 trade data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert two blank rows at row 5 (shifts the old row 5 down to row 7).
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# The inserted rows pick up formatting copied from the row above (the bold
# header row) - wipe that so the new rows start out completely blank, same
# as any other untouched data row.
$ws.Rows.Item(5).Clear()
$ws.Rows.Item(6).Clear()

# --- Row 5: first Feb-09 buy ---
$ws.Range("A5").Value = 46062
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 152.24
$ws.Range("F5").Value = 6120
$ws.Range("G5").Value = "CN#252611665409"
$ws.Range("I5").Value = 30.4
$ws.Range("J5").Formula = '=Index!$C$2'

# --- Row 6: second Feb-09 buy ---
$ws.Range("A6").Value = 46062
$ws.Range("A6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B6").Value = "NSE"
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 149.76
$ws.Range("F6").Value = 5990.4
$ws.Range("G6").Value = "~"
$ws.Range("J6").Formula = '=Index!$C$2'
